# Edit workbook according to the target diff:
#  - workbook.xml: activeTab 2 -> 3 (Sweet Delights becomes active tab)
#  - sheet1 (Cream Cakes): selection C8 -> D13 ; D12 price 500 -> 450
#  - sheet3 (Pastries): no longer the selected tab
#  - sheet4 (Sweet Delights): dimension grows to AC10 ; selection C13 -> D11 ;
#       becomes the selected tab ; D3 price 45 -> 35 ; new rows 6-10 (cookies)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Cream Cakes")
$ws3 = $wb.Worksheets.Item("Pastries")
$ws4 = $wb.Worksheets.Item("Sweet Delights")

# ---------------------------------------------------------------------------
# Cream Cakes: update price in D12 (500 -> 450)
# ---------------------------------------------------------------------------
$ws1.Cells.Item(12, 4).Value = 450

# ---------------------------------------------------------------------------
# Sweet Delights: update price in D3 (45 -> 35)
# ---------------------------------------------------------------------------
$ws4.Cells.Item(3, 4).Value = 35

# ---------------------------------------------------------------------------
# Sweet Delights: add 5 new rows (sd5..sd9 cookies) after existing row 5
# ---------------------------------------------------------------------------

# Copy formatting from row 5 down into the new rows 6-10 (columns A-E and G-H,
# column F is intentionally left untouched/empty just like row 5).
$ws4.Range("A5:E5").Copy() | Out-Null
$ws4.Range("A6:E10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws4.Range("G5:H5").Copy() | Out-Null
$ws4.Range("G6:H10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 6; $r -le 10; $r++) {
    $ws4.Rows.Item($r).RowHeight = 15.75
}

# Fill in the new values, following the same entry order the original
# workbook used (ids/names first, then the shared "For 100 gm" price unit,
# then remaining ids/names, then image paths, then prices, then remaining
# price-units, then stock/discount flags) so shared strings line up.
$ws4.Cells.Item(6, 1).Value = "sd5"
$ws4.Cells.Item(6, 2).Value = "jeera cookies"
$ws4.Cells.Item(6, 5).Value = "For 100 gm"

$ws4.Cells.Item(7, 1).Value = "sd6"
$ws4.Cells.Item(8, 1).Value = "sd7"
$ws4.Cells.Item(9, 1).Value = "sd8"

$ws4.Cells.Item(7, 2).Value = "atta cookies"
$ws4.Cells.Item(8, 2).Value = "chocolate cookies"
$ws4.Cells.Item(9, 2).Value = "coconut cookies"

$ws4.Cells.Item(10, 1).Value = "sd9"
$ws4.Cells.Item(10, 2).Value = "ragi oats cookies"

$ws4.Cells.Item(7, 3).Value  = "sweet-delights/atta-cookies.jpg"
$ws4.Cells.Item(8, 3).Value  = "sweet-delights/chocolate-cookies.jpg"
$ws4.Cells.Item(9, 3).Value  = "sweet-delights/coconut-cookies.jpg"
$ws4.Cells.Item(10, 3).Value = "sweet-delights/ragi-oats-cookies.jpg"
$ws4.Cells.Item(6, 3).Value  = "sweet-delights/jeera-cookies.jpg"

$ws4.Cells.Item(6, 4).Value  = 30
$ws4.Cells.Item(7, 4).Value  = 30
$ws4.Cells.Item(8, 4).Value  = 40
$ws4.Cells.Item(9, 4).Value  = 35
$ws4.Cells.Item(10, 4).Value = 40

$ws4.Cells.Item(7, 5).Value  = "For 100 gm"
$ws4.Cells.Item(8, 5).Value  = "For 100 gm"
$ws4.Cells.Item(9, 5).Value  = "For 100 gm"
$ws4.Cells.Item(10, 5).Value = "For 100 gm"

for ($r = 6; $r -le 10; $r++) {
    $ws4.Cells.Item($r, 7).Value = "yes"
    $ws4.Cells.Item($r, 8).Value = "no"
}

# ---------------------------------------------------------------------------
# Sheet selections / active cells
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D13").Select() | Out-Null

$ws3.Activate()
$ws3.Range("D9").Select() | Out-Null

# Sweet Delights becomes the active / selected tab last.
$ws4.Activate()
$ws4.Range("D11").Select() | Out-Null
